# Delete/Block/Closed - Account scripts - Duplicate
#
# Renames a handful of long account-type labels to their shorter
# "display" forms on the AccountHeirarchy sheet, replaces the stray
# VTA2 / VTA4 placeholders with the canonical "VTA" value, and adds a
# new duplicated VTA row (row 8) mirroring rows 6/7. Finally the
# AccountHeirarchy sheet (instead of Template) becomes the active /
# selected sheet, with A6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountHeirarchy")

# Shorten the account name labels used in column B (Parent) and
# column E (Name of the account). Every occurrence of each long-form
# label is replaced with its short form.
$ws.Range("E3").Value = "ShadowAcc"
$ws.Range("B4").Value = "ShadowAcc"
$ws.Range("B5").Value = "ShadowAcc"

$ws.Range("E2").Value = "Root Acc"
$ws.Range("B3").Value = "Root Acc"

$ws.Range("E4").Value = "Agg One"
$ws.Range("B6").Value = "Agg One"

$ws.Range("E5").Value = "Agg Two"
$ws.Range("B7").Value = "Agg Two"

# Replace the duplicate VTA2 / VTA4 codes with the plain "VTA" value.
$ws.Range("E6").Value = "VTA"
$ws.Range("E7").Value = "VTA"

# New duplicated VTA account row.
$ws.Range("B8").Value = "ShadowAcc"
$ws.Range("C8").Value = "Virtual Transaction Account"
$ws.Range("E8").Value = "VTA"
$ws.Range("F8").Value = "VTA"
$ws.Range("G8").Value = "NORWAY"
$ws.Range("H8").Value = "NOK"
$ws.Range("V8").Style = "Hyperlink"

# Make AccountHeirarchy the active sheet/tab again and select A6
# (previously Template/sheet4 held the tab selection).
$ws.Activate()
$ws.Range("A6").Select()
